$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 5986.5454
$ws.Range("I43").Value = 2616.6667
$ws.Range("J43").Value = 7250.25
$ws.Range("K43").Value = 2616.6667
$ws.Range("L43").Value = 7250.25
$ws.Range("M43").Value = -2547.6667
$ws.Range("N43").Value = -7388.25
$ws.Range("H98").Value = 3788.6428
$ws.Range("I98").Value = 1449
$ws.Range("J98").Value = 8000
$ws.Range("K98").Value = 1449
$ws.Range("L98").Value = 8000
$ws.Range("M98").Value = 49
$ws.Range("N98").Value = -10996
$ws.Range("H122").Value = 3788.6428
$ws.Range("I122").Value = 1449
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 4347
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -1897
$ws.Range("N122").Value = -28900
$ws.Range("H137").Value = 2912.0652
$ws.Range("I137").Value = 2767.7942
$ws.Range("J137").Value = 3320.8333
$ws.Range("K137").Value = 8303.382599999999
$ws.Range("L137").Value = 9962.499899999999
$ws.Range("M137").Value = -5753.382599999999
$ws.Range("N137").Value = -15062.4999

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3664.4285
$ws.Range("I2").Value = 2846.889
$ws.Range("J2").Value = 5136
$ws.Range("K2").Value = 2846.889
$ws.Range("L2").Value = 5136
$ws.Range("M2").Value = -2733.889
$ws.Range("N2").Value = -5362
$ws.Range("H74").Value = 381579.44
$ws.Range("I74").Value = 589524.3
$ws.Range("J74").Value = 86990.914
$ws.Range("K74").Value = 589524.3
$ws.Range("L74").Value = 86990.914
$ws.Range("M74").Value = -588650.3
$ws.Range("N74").Value = -88738.914
$ws.Range("H77").Value = 381579.44
$ws.Range("I77").Value = 589524.3
$ws.Range("J77").Value = 86990.914
$ws.Range("K77").Value = 2947621.5
$ws.Range("L77").Value = 434954.57
$ws.Range("M77").Value = -2943253.5
$ws.Range("N77").Value = -443690.57
$ws.Range("H116").Value = 3664.4285
$ws.Range("I116").Value = 2846.889
$ws.Range("J116").Value = 5136
$ws.Range("K116").Value = 2846.889
$ws.Range("L116").Value = 5136
$ws.Range("M116").Value = -552.8890000000001
$ws.Range("N116").Value = -9724
$ws.Range("H132").Value = 14584.136
$ws.Range("I132").Value = 17438.656
$ws.Range("K132").Value = 52315.96799999999
$ws.Range("M132").Value = -49785.96799999999

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3664.4285
$ws.Range("I3").Value = 2846.889
$ws.Range("J3").Value = 5136
$ws.Range("K3").Value = 2846.889
$ws.Range("L3").Value = 5136
$ws.Range("M3").Value = -2732.889
$ws.Range("N3").Value = -5364
$ws.Range("H134").Value = 4092
$ws.Range("I134").Value = 4228.091
$ws.Range("J134").Value = 3878.1428
$ws.Range("K134").Value = 12684.273
$ws.Range("L134").Value = 11634.4284
$ws.Range("M134").Value = -10149.273
$ws.Range("N134").Value = -16704.4284

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2727.9644
$ws.Range("I132").Value = 1344.8334
$ws.Range("J132").Value = 3765.3125
$ws.Range("K132").Value = 4034.5002
$ws.Range("L132").Value = 11295.9375
$ws.Range("M132").Value = -1504.5002
$ws.Range("N132").Value = -16355.9375

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1579.2
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1579.2
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 4737.6
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -5075.6
$ws.Range("H39").Value = 2950.4
$ws.Range("J39").Value = 2867.3333
$ws.Range("L39").Value = 8601.999899999999
$ws.Range("N39").Value = -9189.999899999999
$ws.Range("H46").Value = 1075.5555
$ws.Range("I46").Value = 936
$ws.Range("J46").Value = 1250
$ws.Range("K46").Value = 2808
$ws.Range("L46").Value = 3750
$ws.Range("M46").Value = -2717
$ws.Range("N46").Value = -3932
$ws.Range("H49").Value = 2551
$ws.Range("J49").Value = 2401.3333
$ws.Range("L49").Value = 7203.999899999999
$ws.Range("N49").Value = -7515.999899999999
$ws.Range("H54").Value = 2563.8
$ws.Range("I54").Value = 2304
$ws.Range("J54").Value = 2628.75
$ws.Range("K54").Value = 6912
$ws.Range("L54").Value = 7886.25
$ws.Range("M54").Value = -6353
$ws.Range("N54").Value = -9004.25
$ws.Range("H55").Value = 3122.2222
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 3122.2222
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 9366.6666
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -9720.6666
$ws.Range("H57").Value = 555.5
$ws.Range("I57").Value = 555.5
$ws.Range("K57").Value = 1666.5
$ws.Range("M57").Value = -1107.5
$ws.Range("H59").Value = 1833.3334
$ws.Range("I59").Value = 500
$ws.Range("J59").Value = 2500
$ws.Range("K59").Value = 1500
$ws.Range("L59").Value = 7500
$ws.Range("M59").Value = -960
$ws.Range("N59").Value = -8580
$ws.Range("H60").Value = 129.5
$ws.Range("I60").Value = 89
$ws.Range("J60").Value = 170
$ws.Range("K60").Value = 267
$ws.Range("L60").Value = 510
$ws.Range("M60").Value = -16
$ws.Range("N60").Value = -1012
$ws.Range("H61").Value = 960
$ws.Range("I61").Value = 400
$ws.Range("J61").Value = 1240
$ws.Range("K61").Value = 1200
$ws.Range("L61").Value = 3720
$ws.Range("M61").Value = -985
$ws.Range("N61").Value = -4150
$ws.Range("H74").Value = 2733.3333
$ws.Range("I74").Value = 200
$ws.Range("K74").Value = 600
$ws.Range("M74").Value = 461
$ws.Range("H77").Value = 2733.3333
$ws.Range("I77").Value = 200
$ws.Range("K77").Value = 1800
$ws.Range("M77").Value = 3504

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8088.8887
$ws.Range("I70").Value = 4005.2632
$ws.Range("J70").Value = 12652.941
$ws.Range("K70").Value = 4005.2632
$ws.Range("L70").Value = 12652.941
$ws.Range("M70").Value = -3735.2632
$ws.Range("N70").Value = -13192.941
$ws.Range("H73").Value = 8088.8887
$ws.Range("I73").Value = 4005.2632
$ws.Range("J73").Value = 12652.941
$ws.Range("K73").Value = 4005.2632
$ws.Range("L73").Value = 12652.941
$ws.Range("M73").Value = -3069.2632
$ws.Range("N73").Value = -14524.941

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2630.077
$ws.Range("I16").Value = 2670.1
$ws.Range("J16").Value = 2496.6667
$ws.Range("K16").Value = 2670.1
$ws.Range("L16").Value = 2496.6667
$ws.Range("M16").Value = -2500.1
$ws.Range("N16").Value = -2836.6667
$ws.Range("H40").Value = 2775
$ws.Range("I40").Value = 2366.6667
$ws.Range("K40").Value = 2366.6667
$ws.Range("M40").Value = -2230.6667
$ws.Range("H61").Value = 1577
$ws.Range("I61").Value = 1102.6666
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 1102.6666
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -900.6666
$ws.Range("N61").Value = -3404
$ws.Range("H113").Value = 1577
$ws.Range("I113").Value = 1102.6666
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 1102.6666
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 1067.3334
$ws.Range("N113").Value = -7340
$ws.Range("H132").Value = 4386.9253
$ws.Range("I132").Value = 1267.4222
$ws.Range("J132").Value = 10767.728
$ws.Range("K132").Value = 3802.2666
$ws.Range("L132").Value = 32303.184
$ws.Range("M132").Value = -1272.2666
$ws.Range("N132").Value = -37363.18399999999

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1466.0625
$ws.Range("I132").Value = 975.525
$ws.Range("J132").Value = 2283.625
$ws.Range("K132").Value = 2926.575
$ws.Range("L132").Value = 6850.875
$ws.Range("M132").Value = -396.5749999999998
$ws.Range("N132").Value = -11910.875
$ws.Range("H135").Value = 31050.385
$ws.Range("J135").Value = 31050.385
$ws.Range("L135").Value = 31050.385
$ws.Range("N135").Value = -41190.38499999999
